# FleetResParam.xlsx – SLEP lines fix
# - NUM_INSTRUCTOR (C) collapses to 40 for every row
# - s_o_c (D) collapses to 25 for rows 7-10 (was 50)
# - SLEPspots (Q), addHours (R), SLEP_or_not (P), Stagger (S) and TTR (T)
#   get refreshed per-row values/formulas
# - two new columns (U = sunDownDate, V = sunDownLength) are populated
# - rows 11-16 (the old ip=75 block) are removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- NUM_INSTRUCTOR (column C) : every remaining row becomes 40 ----
$ws.Range("C2:C10").Value = 40

# ---- s_o_c (column D) : rows 7-10 move from 50 down to 25 ----
$ws.Range("D7:D10").Value = 25

# ---- SLEPspots (column Q) ----
$ws.Range("Q2").Value = 4
$ws.Range("Q3").Value = 4
$ws.Range("Q6").Value = 4
$ws.Range("Q7").Value = 4
$ws.Range("Q10").Value = 4

# ---- addHours (column R) : rows 6-10 flip to FALSE ----
$ws.Range("R6:R10").Value = $false

# ---- SLEP_or_not (column P) : row 10 flips to FALSE ----
$ws.Range("P10").Value = $false

# ---- Stagger (column S) : rows with the shorter window drop to 10800 ----
$ws.Range("S3").Value = 10800
$ws.Range("S5").Value = 10800
$ws.Range("S7").Value = 10800
$ws.Range("S9").Value = 10800
$ws.Range("S10").Value = 10800

# ---- TTR (column T) : re-enter explicit formulas per row ----
$ws.Range("T2").Formula = "=24*30*9"
$ws.Range("T3").Formula = "=24*30*6"
$ws.Range("T4").Formula = "=24*30*9"
$ws.Range("T5").Formula = "=24*30*6"
$ws.Range("T6").Formula = "=24*30*9"
$ws.Range("T7").Formula = "=24*30*6"
$ws.Range("T8").Formula = "=24*30*9"
$ws.Range("T9").Formula = "=24*30*6"
$ws.Range("T10").Formula = "=24*30*6"

# ---- sunDownDate (column U) / sunDownLength (column V) : brand new ----
# Rows 2 and 3 get their own standalone formula; rows 4-10 are entered as one
# fill so the engine records row 4 as the shared-formula anchor (matches how
# the author typed 2-3 by hand, then filled 4:10 down).
$ws.Range("U2").Formula = "=22*365*24"
$ws.Range("V2").Formula = "=5*365*24"
$ws.Range("U3").Formula = "=22*365*24"
$ws.Range("V3").Formula = "=5*365*24"
$ws.Range("U4:U10").Formula = "=22*365*24"
$ws.Range("V4:V10").Formula = "=5*365*24"

# ---- drop the old ip=75 block (rows 11-16) ----
$ws.Rows("11:16").Delete()

# ---- refresh the view: active cell moves to H2, no frozen scroll offset ----
[void]$ws.Range("H2").Select()
